$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292:305 down to 293:306.
$ws.Rows("292:292").Insert()

# Populate the newly inserted row 292 with the new record.
$ws.Range("A292").Value = 10
$ws.Range("B292").Value = "Vega Modelo de Temuco"
$ws.Range("C292").Value = "La Araucanía"
$ws.Range("D292").Value = 44753
$ws.Range("E292").Value = 9
$ws.Range("F292").Value = 100112001
$ws.Range("G292").Value = "Berenjena"
$ws.Range("H292").Value = "Sin especificar"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 180
$ws.Range("K292").Value = 13000
$ws.Range("L292").Value = 14000
$ws.Range("M292").Value = 13556
$ws.Range("N292").Value = "$/caja 60 unidades"
$ws.Range("O292").Value = "Región de Arica y Parinacota"
$ws.Range("P292").Value = 226
$ws.Range("Q292").Value = 60
$ws.Range("R292").Value = "Hortaliza"

# Match the date-number-format style used by the rest of column D.
$ws.Range("D292").NumberFormat = "YYYY-MM-DD HH:MM:SS"
